$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.960.13"
$ws.Range("E2").Value = "  -0.28%  "

$ws.Range("D3").Value = "1.640.04"
$ws.Range("E3").Value = "  +0.17%  "

$ws.Range("E4").Value = "  +0.80%  "

$ws.Range("D5").Value = "'215.01"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("E6").Value = "  +0.64%  "

$ws.Range("E8").Value = "  -0.45%  "

$ws.Range("D9").Value = "'0.0638"
$ws.Range("E9").Value = "  +0.80%  "

$ws.Range("D10").Value = "19.66"
$ws.Range("E10").Value = "  -0.72%  "

$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("B12").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C12").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D12").Value = "1.867.00"
$ws.Range("E12").Value = "  +0.30%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.25"
$ws.Range("E13").Value = "  -0.08%  "

$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.640.47"
$ws.Range("E14").Value = "  +0.54%  "

$ws.Range("E15").Value = "  -1.41%  "

$ws.Range("E16").Value = "  -0.16%  "

$ws.Range("D17").Value = "'62.70"
$ws.Range("E17").Value = "  -0.87%  "

$ws.Range("D18").Value = "25.961.13"
$ws.Range("E18").Value = "  -0.15%  "

$ws.Range("E19").Value = "  +0.86%  "

$ws.Range("D20").Value = "'194.27"
$ws.Range("E20").Value = "  +0.68%  "

$ws.Range("D21").Value = "'4.38"
$ws.Range("E21").Value = "  -1.53%  "

$ws.Range("E22").Value = "  -0.74%  "

$ws.Range("E23").Value = "  -1.18%  "

$ws.Range("D24").Value = "'144.16"
$ws.Range("E24").Value = "  +1.40%  "

$ws.Range("D25").Value = "'1.78"
$ws.Range("E25").Value = "  -0.37%  "

$ws.Range("E26").Value = "  +0.93%  "

$ws.Range("E27").Value = "  +1.99%  "

$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("D29").Value = "'15.50"
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("E30").Value = "  +0.03%  "

$ws.Range("E31").Value = "  +0.50%  "

$ws.Range("E32").Value = "  -1.19%  "

$ws.Range("E33").Value = "  -0.15%  "

$ws.Range("D34").Value = "'1.55"
$ws.Range("E34").Value = "  -2.74%  "

$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("D36").Value = "'0.905"
$ws.Range("E36").Value = "  -0.16%  "

$ws.Range("D37").Value = "1.140.16"
$ws.Range("E37").Value = "  +0.00%  "

$ws.Range("E38").Value = "  -0.87%  "

$ws.Range("D39").Value = "'2.45"
$ws.Range("E39").Value = "  -1.59%  "

$ws.Range("E40").Value = "  +0.24%  "

$ws.Range("D41").Value = "'99.40"
$ws.Range("E41").Value = "  -0.85%  "

$ws.Range("E42").Value = "  +1.23%  "

$ws.Range("E43").Value = "  -2.92%  "

$ws.Range("D44").Value = "1.776.56"
$ws.Range("E44").Value = "  +0.38%  "

$ws.Range("D45").Value = "0.0₆0115"
$ws.Range("E45").Value = "  +8.01%  "

$ws.Range("E46").Value = "  +1.35%  "

$ws.Range("E47").Value = "  +2.63%  "

$ws.Range("E48").Value = "  -0.51%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'7.67"
$ws.Range("E49").Value = "  +0.13%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "'0.415"
$ws.Range("E50").Value = "  -0.28%  "

$ws.Range("D51").Value = "'0.0964"
$ws.Range("E51").Value = "  -0.86%  "

